# Auto-generated Excel COM-interop script
# Applies targeted numeric cell updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 8022455.5
$ws.Range("J112").Value = 8022455.5
$ws.Range("L112").Value = 24067366.5
$ws.Range("N112").Value = -24069582.5
$ws.Range("H138").Value = 5776739
$ws.Range("J138").Value = 8477133
$ws.Range("L138").Value = 25431399
$ws.Range("N138").Value = -25441679
$ws.Range("H141").Value = 2865.5
$ws.Range("I141").Value = 2728.6
$ws.Range("J141").Value = 3550
$ws.Range("K141").Value = 8185.799999999999
$ws.Range("L141").Value = 10650
$ws.Range("M141").Value = -3005.799999999999
$ws.Range("N141").Value = -21010

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3785.9
$ws.Range("I2").Value = 4920.2593
$ws.Range("J2").Value = 1429.9231
$ws.Range("K2").Value = 4920.2593
$ws.Range("L2").Value = 1429.9231
$ws.Range("M2").Value = -4807.2593
$ws.Range("N2").Value = -1655.9231
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H97").Value = 18519248
$ws.Range("I97").Value = 30303824
$ws.Range("J97").Value = 629.8570999999999
$ws.Range("K97").Value = 30303824
$ws.Range("L97").Value = 629.8570999999999
$ws.Range("M97").Value = -30303328
$ws.Range("N97").Value = -1621.8571
$ws.Range("H102").Value = 1545.2727
$ws.Range("I102").Value = 1242.8572
$ws.Range("J102").Value = 2074.5
$ws.Range("K102").Value = 1242.8572
$ws.Range("L102").Value = 2074.5
$ws.Range("M102").Value = 379.1428000000001
$ws.Range("N102").Value = -5318.5
$ws.Range("H107").Value = 20228
$ws.Range("J107").Value = 20228
$ws.Range("L107").Value = 20228
$ws.Range("N107").Value = -27908
$ws.Range("H109").Value = 130075.4
$ws.Range("J109").Value = 130075.4
$ws.Range("L109").Value = 130075.4
$ws.Range("N109").Value = -132849.4
$ws.Range("H116").Value = 3785.9
$ws.Range("I116").Value = 4920.2593
$ws.Range("J116").Value = 1429.9231
$ws.Range("K116").Value = 4920.2593
$ws.Range("L116").Value = 1429.9231
$ws.Range("M116").Value = -2626.2593
$ws.Range("N116").Value = -6017.9231
$ws.Range("H133").Value = 52662.125
$ws.Range("J133").Value = 52662.125
$ws.Range("L133").Value = 52662.125
$ws.Range("N133").Value = -57722.125
$ws.Range("H139").Value = 65488.43
$ws.Range("J139").Value = 65488.43
$ws.Range("L139").Value = 65488.43
$ws.Range("N139").Value = -75768.42999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3785.9
$ws.Range("I3").Value = 4920.2593
$ws.Range("J3").Value = 1429.9231
$ws.Range("K3").Value = 4920.2593
$ws.Range("L3").Value = 1429.9231
$ws.Range("M3").Value = -4806.2593
$ws.Range("N3").Value = -1657.9231
$ws.Range("H86").Value = 9979
$ws.Range("I86").Value = 1902.5
$ws.Range("K86").Value = 1902.5
$ws.Range("M86").Value = -779.5
$ws.Range("H89").Value = 9979
$ws.Range("I89").Value = 1902.5
$ws.Range("K89").Value = 9512.5
$ws.Range("M89").Value = -3896.5
$ws.Range("H94").Value = 945.9474
$ws.Range("I94").Value = 989.5714
$ws.Range("J94").Value = 823.8
$ws.Range("K94").Value = 989.5714
$ws.Range("L94").Value = 823.8
$ws.Range("M94").Value = -538.5714
$ws.Range("N94").Value = -1725.8
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 49000
$ws.Range("L133").Value = 49000
$ws.Range("N133").Value = -59120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1609.0741
$ws.Range("I31").Value = 975.1667
$ws.Range("J31").Value = 2876.889
$ws.Range("K31").Value = 975.1667
$ws.Range("L31").Value = 2876.889
$ws.Range("M31").Value = -680.1667
$ws.Range("N31").Value = -3466.889
$ws.Range("H34").Value = 1609.0741
$ws.Range("I34").Value = 975.1667
$ws.Range("J34").Value = 2876.889
$ws.Range("K34").Value = 975.1667
$ws.Range("L34").Value = 2876.889
$ws.Range("M34").Value = -773.1667
$ws.Range("N34").Value = -3280.889
$ws.Range("H62").Value = 26944.445
$ws.Range("I62").Value = 42100
$ws.Range("K62").Value = 42100
$ws.Range("M62").Value = -41476
$ws.Range("H65").Value = 26944.445
$ws.Range("I65").Value = 42100
$ws.Range("K65").Value = 210500
$ws.Range("M65").Value = -207380
$ws.Range("H132").Value = 4411.769
$ws.Range("I132").Value = 3464.8
$ws.Range("K132").Value = 10394.4
$ws.Range("M132").Value = -7864.400000000001
$ws.Range("H134").Value = 3068.739
$ws.Range("I134").Value = 1370.8667
$ws.Range("J134").Value = 6252.25
$ws.Range("K134").Value = 4112.6001
$ws.Range("L134").Value = 18756.75
$ws.Range("M134").Value = -1577.6001
$ws.Range("N134").Value = -23826.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 836.6
$ws.Range("I97").Value = 795.75
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2387.25
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1891.25
$ws.Range("N97").Value = -3992
$ws.Range("H122").Value = 962.0952
$ws.Range("I122").Value = 785.3
$ws.Range("J122").Value = 1122.8182
$ws.Range("K122").Value = 7067.7
$ws.Range("L122").Value = 10105.3638
$ws.Range("M122").Value = -4617.7
$ws.Range("N122").Value = -15005.3638
$ws.Range("H131").Value = 3191.5264
$ws.Range("I131").Value = 783.3333
$ws.Range("J131").Value = 3474.843
$ws.Range("K131").Value = 2349.9999
$ws.Range("L131").Value = 10424.529
$ws.Range("M131").Value = 2690.0001
$ws.Range("N131").Value = -20504.529
$ws.Range("H137").Value = 6318560
$ws.Range("I137").Value = 12505803
$ws.Range("J137").Value = 131317
$ws.Range("K137").Value = 37517409
$ws.Range("L137").Value = 393951
$ws.Range("M137").Value = -37512309
$ws.Range("N137").Value = -404151
$ws.Range("H139").Value = 35716616
$ws.Range("I139").Value = 38463470
$ws.Range("J139").Value = 7500
$ws.Range("K139").Value = 115390410
$ws.Range("L139").Value = 22500
$ws.Range("M139").Value = -115385270
$ws.Range("N139").Value = -32780
$ws.Range("H141").Value = 8642
$ws.Range("I141").Value = 9570
$ws.Range("J141").Value = 7250
$ws.Range("K141").Value = 28710
$ws.Range("L141").Value = 21750
$ws.Range("M141").Value = -23530
$ws.Range("N141").Value = -32110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1057.6364
$ws.Range("I107").Value = 1083.4
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1083.4
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 836.5999999999999
$ws.Range("N107").Value = -4640
$ws.Range("H123").Value = 10321.267
$ws.Range("J123").Value = 10321.267
$ws.Range("L123").Value = 10321.267
$ws.Range("N123").Value = -15221.267
$ws.Range("H126").Value = 2010.7742
$ws.Range("I126").Value = 1408.3636
$ws.Range("J126").Value = 2342.1
$ws.Range("K126").Value = 4225.0908
$ws.Range("L126").Value = 7026.299999999999
$ws.Range("M126").Value = -1755.0908
$ws.Range("N126").Value = -11966.3
$ws.Range("H132").Value = 2955.3257
$ws.Range("I132").Value = 2378.913
$ws.Range("J132").Value = 3618.2
$ws.Range("K132").Value = 7136.739
$ws.Range("L132").Value = 10854.6
$ws.Range("M132").Value = -4606.739
$ws.Range("N132").Value = -15914.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3535.7144
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3535.7144
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3535.7144
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -3759.7144
$ws.Range("H40").Value = 4503.125
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4503.125
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4503.125
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -4775.125
$ws.Range("H126").Value = 3535.7144
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3535.7144
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10607.1432
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -15547.1432
$ws.Range("H133").Value = 48066
$ws.Range("J133").Value = 48066
$ws.Range("L133").Value = 48066
$ws.Range("N133").Value = -53126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1219.4546
$ws.Range("I113").Value = 1452.3334
$ws.Range("J113").Value = 171.5
$ws.Range("K113").Value = 4357.0002
$ws.Range("L113").Value = 514.5
$ws.Range("M113").Value = -2187.0002
$ws.Range("N113").Value = -4854.5
$ws.Range("H122").Value = 67880.2
$ws.Range("I122").Value = 143858.28
$ws.Range("J122").Value = 1399.375
$ws.Range("K122").Value = 431574.84
$ws.Range("L122").Value = 4198.125
$ws.Range("M122").Value = -429124.84
$ws.Range("N122").Value = -9098.125
$ws.Range("H132").Value = 50004650
$ws.Range("I132").Value = 83337770
$ws.Range("K132").Value = 250013310
$ws.Range("M132").Value = -250010780
